$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.253.59'
$ws.Range('E2').Value = '  -2.36%  '
$ws.Range('D3').Value = '1.869.32'
$ws.Range('E3').Value = '  -1.88%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '318.55'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.80%  '
$ws.Range('E6').Value = '  -0.02%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4381'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -4.58%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3692'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -3.41%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07505'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.63%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.9374'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -4.15%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '21.47'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.66%  '
$ws.Range('D12').Value = '1.888.50'
$ws.Range('E12').Value = '  -1.36%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.722'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -3.00%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.449'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -3.65%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.06862'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.55%  '
$ws.Range('E16').Value = '  -0.11%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '82.05'
$ws.Range('D17').Style = 'Normal'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000009052'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -4.18%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.001'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.09%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '15.95'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -3.96%  '
$ws.Range('D21').Value = '28.239.69'
$ws.Range('E21').Value = '  -2.42%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.129'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -3.19%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.81'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.60%  '
$ws.Range('D24').Value = '2.117.74'
$ws.Range('E24').Value = '  -1.86%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.027'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -3.27%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '154.90'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.03%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.38'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -3.60%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '5.318'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -5.90%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '114.00'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -2.92%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.726'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -6.12%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.09039'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.45%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.7993'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -7.37%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.838'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -4.75%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.173'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -5.76%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.930'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.64%  '
$ws.Range('E36').Value = '  -0.03%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.121'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -2.17%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.05451'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -4.50%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01955'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -4.21%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.915'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +5.14%  '
$ws.Range('B41').Value = 'FraxShare'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '7.110'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -3.72%  '
$ws.Range('B42').Value = 'TheSandbox'
$ws.Range('C42').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.5257'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -4.25%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.1681'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -4.20%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.764'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -5.65%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.06752'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.83%  '
$ws.Range('E46').Value = '  -5.57%  '
$ws.Range('B47').Value = 'RenderToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.987'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -3.30%  '
$ws.Range('B48').Value = 'Quant'
$ws.Range('C48').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '107.86'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.18%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '10.51'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -6.47%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.000002443'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -4.65%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.680'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -5.02%  '
